$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3300.25
$ws.Range("I40").Value = 1666.8334
$ws.Range("J40").Value = 4933.6665
$ws.Range("K40").Value = 1666.8334
$ws.Range("L40").Value = 4933.6665
$ws.Range("M40").Value = -1491.8334
$ws.Range("N40").Value = -5283.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8093.5933
$ws.Range("I32").Value = 2655.9714
$ws.Range("J32").Value = 26219
$ws.Range("K32").Value = 2655.9714
$ws.Range("L32").Value = 26219
$ws.Range("M32").Value = -2368.9714
$ws.Range("N32").Value = -26793

$ws.Range("H88").Value = 166668740
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 166668740
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 166668740
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -166669552

$ws.Range("H91").Value = 166668740
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 166668740
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 166668740
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -166671548

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 3501
$ws.Range("J15").Value = 3501
$ws.Range("L15").Value = 3501
$ws.Range("N15").Value = -3955

$ws.Range("H54").Value = 11538.667
$ws.Range("I54").Value = 789.6
$ws.Range("J54").Value = 24975
$ws.Range("K54").Value = 789.6
$ws.Range("L54").Value = 24975
$ws.Range("M54").Value = -305.6
$ws.Range("N54").Value = -25943

$ws.Range("H82").Value = 14042.818
$ws.Range("J82").Value = 32300
$ws.Range("L82").Value = 32300
$ws.Range("N82").Value = -33066

$ws.Range("H85").Value = 14042.818
$ws.Range("J85").Value = 32300
$ws.Range("L85").Value = 32300
$ws.Range("N85").Value = -34952

$ws.Range("H86").Value = 25003500
$ws.Range("I86").Value = 28574000
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 28574000
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -28572877
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 25003500
$ws.Range("I89").Value = 28574000
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 142870000
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -142864384
$ws.Range("N89").Value = -61232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 264806.78
$ws.Range("I86").Value = 455982.47
$ws.Range("J86").Value = 1940.25
$ws.Range("K86").Value = 455982.47
$ws.Range("L86").Value = 1940.25
$ws.Range("M86").Value = -454859.47
$ws.Range("N86").Value = -4186.25

$ws.Range("H89").Value = 264806.78
$ws.Range("I89").Value = 455982.47
$ws.Range("J89").Value = 1940.25
$ws.Range("K89").Value = 2279912.35
$ws.Range("L89").Value = 9701.25
$ws.Range("M89").Value = -2274296.35
$ws.Range("N89").Value = -20933.25

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 32263196
$ws.Range("I58").Value = 2002.5
$ws.Range("J58").Value = 34488108
$ws.Range("K58").Value = 6007.5
$ws.Range("L58").Value = 103464324
$ws.Range("M58").Value = -5879.5
$ws.Range("N58").Value = -103464580

$ws.Range("H64").Value = 52636092
$ws.Range("I64").Value = 2909
$ws.Range("J64").Value = 66671610
$ws.Range("K64").Value = 8727
$ws.Range("L64").Value = 200014830
$ws.Range("M64").Value = -8457
$ws.Range("N64").Value = -200015370

$ws.Range("H67").Value = 52636092
$ws.Range("I67").Value = 2909
$ws.Range("J67").Value = 66671610
$ws.Range("K67").Value = 8727
$ws.Range("L67").Value = 200014830
$ws.Range("M67").Value = -7791
$ws.Range("N67").Value = -200016702

$ws.Range("H68").Value = 1068.9175
$ws.Range("I68").Value = 597.4909
$ws.Range("J68").Value = 1686.262
$ws.Range("K68").Value = 1792.4727
$ws.Range("L68").Value = 5058.786
$ws.Range("M68").Value = -981.4727
$ws.Range("N68").Value = -6680.786

$ws.Range("H70").Value = 4141.05
$ws.Range("I70").Value = 965.125
$ws.Range("J70").Value = 6258.3335
$ws.Range("K70").Value = 2895.375
$ws.Range("L70").Value = 18775.0005
$ws.Range("M70").Value = -2580.375
$ws.Range("N70").Value = -19405.0005

$ws.Range("H71").Value = 1068.9175
$ws.Range("I71").Value = 597.4909
$ws.Range("J71").Value = 1686.262
$ws.Range("K71").Value = 5377.4181
$ws.Range("L71").Value = 15176.358
$ws.Range("M71").Value = -1321.4181
$ws.Range("N71").Value = -23288.358

$ws.Range("H73").Value = 4141.05
$ws.Range("I73").Value = 965.125
$ws.Range("J73").Value = 6258.3335
$ws.Range("K73").Value = 2895.375
$ws.Range("L73").Value = 18775.0005
$ws.Range("M73").Value = -1803.375
$ws.Range("N73").Value = -20959.0005

$ws.Range("H76").Value = 7314.2856
$ws.Range("I76").Value = 6800
$ws.Range("J76").Value = 7400
$ws.Range("K76").Value = 20400
$ws.Range("L76").Value = 22200
$ws.Range("M76").Value = -20017
$ws.Range("N76").Value = -22966

$ws.Range("H79").Value = 7314.2856
$ws.Range("I79").Value = 6800
$ws.Range("J79").Value = 7400
$ws.Range("K79").Value = 20400
$ws.Range("L79").Value = 22200
$ws.Range("M79").Value = -19074
$ws.Range("N79").Value = -24852

$ws.Range("H102").Value = 3500
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1956.5454
$ws.Range("I122").Value = 2146
$ws.Range("J122").Value = 1104
$ws.Range("K122").Value = 6438
$ws.Range("L122").Value = 3312
$ws.Range("M122").Value = -3988
$ws.Range("N122").Value = -8212

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1254.8
$ws.Range("I16").Value = 1283.1666
$ws.Range("J16").Value = 999.5
$ws.Range("K16").Value = 1283.1666
$ws.Range("L16").Value = 999.5
$ws.Range("M16").Value = -1113.1666
$ws.Range("N16").Value = -1339.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 58825396
$ws.Range("I81").Value = 142858640
$ws.Range("J81").Value = 2120.1
$ws.Range("K81").Value = 285717280
$ws.Range("L81").Value = 4240.2
$ws.Range("M81").Value = -285716219
$ws.Range("N81").Value = -6362.2

$ws.Range("H84").Value = 58825396
$ws.Range("I84").Value = 142858640
$ws.Range("J84").Value = 2120.1
$ws.Range("K84").Value = 1428586400
$ws.Range("L84").Value = 21201
$ws.Range("M84").Value = -1428581096
$ws.Range("N84").Value = -31809

$ws.Range("H133").Value = 39875
$ws.Range("J133").Value = 39875
$ws.Range("L133").Value = 39875
$ws.Range("N133").Value = -49995
